# Fix bugs in VAR and ARIMA result sheets: correct fitted/predicted values and
# re-align the Year index (VAR series off-by-one; ARIMA series had an extra leading row).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: y_fitted_on_begin_2016 -----------------------------------
# Years re-aligned (each row now one year earlier) and y_value replaced with the
# corrected fitted series; a new row for 2016 is appended (A1:D36 -> A1:D37).
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$arrS1 = New-Object 'object[,]' 36,2
$arrS1[0,0] = 1981
$arrS1[0,1] = 60.95365071456347
$arrS1[1,0] = 1982
$arrS1[1,1] = 61.8195443180603
$arrS1[2,0] = 1983
$arrS1[2,1] = 62.24993141444148
$arrS1[3,0] = 1984
$arrS1[3,1] = 61.26987512033323
$arrS1[4,0] = 1985
$arrS1[4,1] = 61.12047222365608
$arrS1[5,0] = 1986
$arrS1[5,1] = 61.97603163665062
$arrS1[6,0] = 1987
$arrS1[6,1] = 61.66353022464727
$arrS1[7,0] = 1988
$arrS1[7,1] = 61.97972607560371
$arrS1[8,0] = 1989
$arrS1[8,1] = 61.64512206397297
$arrS1[9,0] = 1990
$arrS1[9,1] = 62.09122112212628
$arrS1[10,0] = 1991
$arrS1[10,1] = 62.89709846252617
$arrS1[11,0] = 1992
$arrS1[11,1] = 62.25863325297033
$arrS1[12,0] = 1993
$arrS1[12,1] = 62.27135693704287
$arrS1[13,0] = 1994
$arrS1[13,1] = 61.99145082286896
$arrS1[14,0] = 1995
$arrS1[14,1] = 62.01860559951943
$arrS1[15,0] = 1996
$arrS1[15,1] = 61.93121223755041
$arrS1[16,0] = 1997
$arrS1[16,1] = 62.41476452945855
$arrS1[17,0] = 1998
$arrS1[17,1] = 62.71028554784527
$arrS1[18,0] = 1999
$arrS1[18,1] = 63.03850128409395
$arrS1[19,0] = 2000
$arrS1[19,1] = 62.79951243648566
$arrS1[20,0] = 2001
$arrS1[20,1] = 62.88720378517831
$arrS1[21,0] = 2002
$arrS1[21,1] = 64.31485821230103
$arrS1[22,0] = 2003
$arrS1[22,1] = 64.51822710068943
$arrS1[23,0] = 2004
$arrS1[23,1] = 64.21623994776404
$arrS1[24,0] = 2005
$arrS1[24,1] = 64.46882645137748
$arrS1[25,0] = 2006
$arrS1[25,1] = 65.40636586527529
$arrS1[26,0] = 2007
$arrS1[26,1] = 65.95471166995229
$arrS1[27,0] = 2008
$arrS1[27,1] = 67.34865122902352
$arrS1[28,0] = 2009
$arrS1[28,1] = 69.04289359310661
$arrS1[29,0] = 2010
$arrS1[29,1] = 68.68739601554691
$arrS1[30,0] = 2011
$arrS1[30,1] = 68.54078951276723
$arrS1[31,0] = 2012
$arrS1[31,1] = 69.41404839704032
$arrS1[32,0] = 2013
$arrS1[32,1] = 69.72868080977396
$arrS1[33,0] = 2014
$arrS1[33,1] = 70.38165332527123
$arrS1[34,0] = 2015
$arrS1[34,1] = 69.92179580068394
$arrS1[35,0] = 2016
$arrS1[35,1] = 69.82630350559577
$ws.Range("A2:B37").Value2 = $arrS1

# --- Sheet 2: y_pred_on_2017_2021 ---------------------------------------
# Years unchanged; y_value replaced with the corrected predicted series.
$ws = $wb.Worksheets.Item("y_pred_on_2017_2021")
$arrS2 = New-Object 'object[,]' 5,1
$arrS2[0,0] = 70.317884850025
$arrS2[1,0] = 70.29983528181613
$arrS2[2,0] = 70.42744403816261
$arrS2[3,0] = 70.64010587333817
$arrS2[4,0] = 70.88925799956439
$ws.Range("B2:B6").Value2 = $arrS2

# --- Sheet 3: y_fitted_on_begin_2021 ------------------------------------
# Drop the stray leading 1980 row, re-align years, and replace y_value with the
# corrected fitted series (A1:D43 -> A1:D42).
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Rows.Item(43).Delete()
$arrS3 = New-Object 'object[,]' 41,2
$arrS3[0,0] = 1981
$arrS3[0,1] = 61.03458975610459
$arrS3[1,0] = 1982
$arrS3[1,1] = 61.93727853491535
$arrS3[2,0] = 1983
$arrS3[2,1] = 62.29648723528988
$arrS3[3,0] = 1984
$arrS3[3,1] = 61.28645049894657
$arrS3[4,0] = 1985
$arrS3[4,1] = 61.1707887784693
$arrS3[5,0] = 1986
$arrS3[5,1] = 61.86306153309339
$arrS3[6,0] = 1987
$arrS3[6,1] = 61.5698312519966
$arrS3[7,0] = 1988
$arrS3[7,1] = 61.91065119412164
$arrS3[8,0] = 1989
$arrS3[8,1] = 61.64044896513014
$arrS3[9,0] = 1990
$arrS3[9,1] = 62.04357613692963
$arrS3[10,0] = 1991
$arrS3[10,1] = 62.80782942276345
$arrS3[11,0] = 1992
$arrS3[11,1] = 62.18863323403155
$arrS3[12,0] = 1993
$arrS3[12,1] = 62.20118406430418
$arrS3[13,0] = 1994
$arrS3[13,1] = 61.94060784045826
$arrS3[14,0] = 1995
$arrS3[14,1] = 61.99788548399147
$arrS3[15,0] = 1996
$arrS3[15,1] = 61.92393327572724
$arrS3[16,0] = 1997
$arrS3[16,1] = 62.41200312521156
$arrS3[17,0] = 1998
$arrS3[17,1] = 62.70277123528523
$arrS3[18,0] = 1999
$arrS3[18,1] = 63.01387683137063
$arrS3[19,0] = 2000
$arrS3[19,1] = 62.77288123908905
$arrS3[20,0] = 2001
$arrS3[20,1] = 63.0119222186326
$arrS3[21,0] = 2002
$arrS3[21,1] = 64.36519091027596
$arrS3[22,0] = 2003
$arrS3[22,1] = 64.5045121411609
$arrS3[23,0] = 2004
$arrS3[23,1] = 64.18285550698647
$arrS3[24,0] = 2005
$arrS3[24,1] = 64.5664099915215
$arrS3[25,0] = 2006
$arrS3[25,1] = 65.48630881746756
$arrS3[26,0] = 2007
$arrS3[26,1] = 66.0630589784264
$arrS3[27,0] = 2008
$arrS3[27,1] = 67.48430624661624
$arrS3[28,0] = 2009
$arrS3[28,1] = 69.07672587453733
$arrS3[29,0] = 2010
$arrS3[29,1] = 68.71720311259307
$arrS3[30,0] = 2011
$arrS3[30,1] = 68.59949510462452
$arrS3[31,0] = 2012
$arrS3[31,1] = 69.44301761567756
$arrS3[32,0] = 2013
$arrS3[32,1] = 69.77783638719501
$arrS3[33,0] = 2014
$arrS3[33,1] = 70.38107879315093
$arrS3[34,0] = 2015
$arrS3[34,1] = 69.91540894406862
$arrS3[35,0] = 2016
$arrS3[35,1] = 69.8429840469887
$arrS3[36,0] = 2017
$arrS3[36,1] = 70.32136847421904
$arrS3[37,0] = 2018
$arrS3[37,1] = 69.83199500032224
$arrS3[38,0] = 2019
$arrS3[38,1] = 69.8444964740936
$arrS3[39,0] = 2020
$arrS3[39,1] = 70.76805854583635
$arrS3[40,0] = 2021
$arrS3[40,1] = 70.47009508756929
$ws.Range("A2:B42").Value2 = $arrS3

# --- Sheet 4: y_pred_on_2022_2026 ---------------------------------------
# Years unchanged; y_value replaced with the corrected predicted series.
$ws = $wb.Worksheets.Item("y_pred_on_2022_2026")
$arrS4 = New-Object 'object[,]' 5,1
$arrS4[0,0] = 68.5981678986218
$arrS4[1,0] = 67.97879844288904
$arrS4[2,0] = 67.84129624785861
$arrS4[3,0] = 67.72691535324546
$arrS4[4,0] = 67.43374798894465
$ws.Range("B2:B6").Value2 = $arrS4
